$d = $word.ActiveDocument

$replacements = @(
    @{old="409÷6="; new="112÷2="},
    @{old="310÷7="; new="649÷8="},
    @{old="608÷4="; new="747÷6="},
    @{old="588÷4="; new="629÷9="},
    @{old="897÷6="; new="886÷4="},
    @{old="859÷7="; new="217÷9="},
    @{old="754÷2="; new="510÷7="},
    @{old="606÷3="; new="785÷2="},
    @{old="655÷4="; new="673÷3="},
    @{old="684÷2="; new="769÷2="},
    @{old="544÷2="; new="529÷4="},
    @{old="438÷2="; new="531÷3="},
    @{old="996÷8="; new="165÷9="},
    @{old="358÷5="; new="440÷6="},
    @{old="986÷7="; new="400÷9="},
    @{old="303÷3="; new="743÷8="},
    @{old="209÷7="; new="793÷6="},
    @{old="411÷8="; new="311÷2="},
    @{old="888÷4="; new="621÷3="},
    @{old="185÷9="; new="540÷9="},
    @{old="402÷4="; new="453÷6="},
    @{old="290÷4="; new="393÷8="},
    @{old="414÷2="; new="477÷6="},
    @{old="844÷5="; new="317÷5="},
    @{old="192÷8="; new="227÷5="}
)

foreach ($r in $replacements) {
    $d.Content.Find.Execute($r.old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $r.new, 2)
}
